$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row at position 35 (shifts old rows 35-36 down to 36-37)
$ws.Rows.Item(35).Insert()

# New row 35: IP in Germany (Dusseldorf) - reuses existing "Dusseldorf, Germany" text
$ws.Range("A35").Value = "93.186.202.39"
$ws.Range("B35").Value = "Dusseldorf, Germany"

# Append new row 38: IP in Canada (Montreal)
$ws.Range("A38").Value = "158.69.158.67"
$ws.Range("B38").Value = "Montreal, Canada"

# Append new row 39: IP in Texas (Round Rock)
$ws.Range("A39").Value = "23.247.148.238"
$ws.Range("B39").Value = "Round Rock, TX"

# Update the selection to match the target state (B39)
$ws.Range("B39").Select()
